$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit the new definitions text (closest reachable width
# to 60.1796875 given this engine's whole-pixel column-width quantization)
$ws.Columns.Item(2).ColumnWidth = 59.25

# Add the "Descripción" definitions in column B for the corresponding rows
$ws.Range("B2").Value  = "ID de orden de compra"
$ws.Range("B3").Value  = "ecommerce a través del cual se realizó la compra"
$ws.Range("B4").Value  = "Tienda en la que se realizó la compra"
$ws.Range("B5").Value  = "Fecha y hora de la transacción"
$ws.Range("B8").Value  = "Categoría del artículo comprado"
$ws.Range("B9").Value  = "Marca del artículo comprado"
$ws.Range("B10").Value = "ID del artículo"
$ws.Range("B11").Value = "Talle del artículo"
$ws.Range("B12").Value = "Nombre del artículo"
$ws.Range("B13").Value = "Cantidad comprada"
$ws.Range("B14").Value = "Precio que pagó el cliente por el artículo"
$ws.Range("B15").Value = "Costo del producto para la empresa"
$ws.Range("B17").Value = "Costo del producto para la empresa (debería ser igual a PrecioCosto)"
$ws.Range("B18").Value = "Valor de venta sin promociones ni descuentos"
$ws.Range("B19").Value = "Forma de envío"
$ws.Range("B22").Value = "Color del artículo"
$ws.Range("B24").Value = "Temática del artículo"
$ws.Range("B25").Value = "Género del artículo"
$ws.Range("B30").Value = "email del comprador (encriptado)"
$ws.Range("B32").Value = "Latitud"
$ws.Range("B33").Value = "Longitud"
$ws.Range("B34").Value = "Método de pago"

# Update the view: scroll so row 16 is at the top-left, and select B32
$ws.Range("B32").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
